$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new record right before current row 99 -------------------------
# Shift existing rows 99.. down by one; the new blank row becomes row 99 and
# the old row 99 becomes row 100.
$ws.Rows(99).Insert()

# Copy the (now shifted) row 100 values into the new row 99, then overwrite
# the two fields that actually differ for the new record (Fecha / Volumen).
$ws.Range("A100:R100").Copy()
$ws.Range("A99").PasteSpecial()
$ws.Range("D99").Value2 = 44567
$ws.Range("J99").Value2 = 80

# --- Insert a second new record right before (the now shifted) row 151 -------
# After the first insert, the record that used to be row 150 now lives at
# row 151. Insert another blank row above it so it moves to row 152, and the
# new blank row becomes row 151.
$ws.Rows(151).Insert()

$ws.Range("A152:R152").Copy()
$ws.Range("A151").PasteSpecial()
$ws.Range("D151").Value2 = 44568
$ws.Range("J151").Value2 = 220

$excel.CutCopyMode = $false
